$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 179 (pushes old rows 179-220 down to 181-222).
$ws.Rows.Item(179).Insert()
$ws.Rows.Item(179).Insert()

$newDate = Get-Date -Year 2022 -Month 8 -Day 25 -Hour 0 -Minute 0 -Second 0

# New row 179: Primera
$ws.Range("A179").Value = 11
$ws.Range("B179").Value = "Vega Monumental Concepción"
$ws.Range("C179").Value = "Bíobío"
$ws.Range("D179").Value = $newDate
$ws.Range("E179").Value = 8
$ws.Range("F179").Value = 100112040
$ws.Range("G179").Value = "Cilantro"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 200
$ws.Range("K179").Value = 700
$ws.Range("L179").Value = 800
$ws.Range("M179").Value = 750
$ws.Range("N179").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O179").Value = "Región de Ñuble"
$ws.Range("P179").Value = 750
$ws.Range("Q179").Value = 1
$ws.Range("R179").Value = "Hortaliza"

# New row 180: Segunda
$ws.Range("A180").Value = 11
$ws.Range("B180").Value = "Vega Monumental Concepción"
$ws.Range("C180").Value = "Bíobío"
$ws.Range("D180").Value = $newDate
$ws.Range("E180").Value = 8
$ws.Range("F180").Value = 100112040
$ws.Range("G180").Value = "Cilantro"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Segunda"
$ws.Range("J180").Value = 100
$ws.Range("K180").Value = 600
$ws.Range("L180").Value = 600
$ws.Range("M180").Value = 600
$ws.Range("N180").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O180").Value = "Región de Ñuble"
$ws.Range("P180").Value = 600
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"
